## Adds two new "postcondition" bullet-style paragraphs (matching the
## existing left-bordered paragraph formatting used in that table cell)
## right after the "All connected alarms have been triggered." paragraph
## in the Postcondition cell of the last table (HandleBreakin use case).
##
## Each new paragraph is built as a small WordprocessingML fragment and
## inserted via Range.InsertXML at a collapsed range, so the formatting
## (pBdr left border + iCs/szCs run-properties, no explicit style) matches
## the sibling paragraphs exactly -- InsertParagraphAfter()/new Range.Text
## in this cell picks up unrelated "Removeable Text" placeholder formatting,
## so we avoid it.

$d = $word.ActiveDocument

$wordMlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-BorderedParagraphXml([string]$text) {
    return '<w:p ' + $wordMlNs + '>' +
        '<w:pPr>' +
            '<w:pBdr><w:left w:val="single" w:sz="18" w:space="4" w:color="auto"/></w:pBdr>' +
            '<w:rPr><w:iCs/><w:szCs w:val="22"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r><w:rPr><w:iCs/><w:szCs w:val="22"/></w:rPr><w:t>' + $text + '</w:t></w:r>' +
    '</w:p>'
}

# --- insert "Phone call has been made." right after the anchor paragraph ---
$anchorText = "All connected alarms have been triggered."
$r1 = $d.Content
$found1 = $r1.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    Write-Output "ERROR: could not find anchor paragraph '$anchorText'"
} else {
    $r1.Collapse(0)   # wdCollapseEnd: collapse to the end of the found range
    $r1.InsertXML((New-BorderedParagraphXml "Phone call has been made.")) | Out-Null

    # --- insert "Event has been logged." right after that new paragraph ---
    $r2 = $d.Content
    $found2 = $r2.Find.Execute("Phone call has been made.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found2) {
        Write-Output "ERROR: could not find the just-inserted 'Phone call has been made.' paragraph"
    } else {
        $r2.Collapse(0)
        $r2.InsertXML((New-BorderedParagraphXml "Event has been logged.")) | Out-Null
        Write-Output "Inserted 'Phone call has been made.' and 'Event has been logged.' paragraphs."
    }
}
